# Update database and roll the reporting window forward by one year:
#   - drop the oldest period (1396/12) and its publish date
#   - shift every remaining period/date one column to the left (D<-E<-F<-G<-H)
#   - populate the newly opened H column (1401/12, published 1402-02-29) with
#     the latest filed figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8: period headers (D:H) ----
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---- Row 9: publish dates (D:H) ----
$ws.Range("D9").Value = "1399-02-22 (14)"
$ws.Range("E9").Value = "1400-03-02 (12)"
$ws.Range("F9").Value = "1401-02-26 (9)"
$ws.Range("G9").Value = "1402-02-29 (8)"
$ws.Range("H9").Value = "1402-02-29"

# ---- Data rows 12-58: D,E,F,G,H values (row, D, E, F, G, H) ----
$rows = @(
    @(12, 128815, 212052, 194476, 767257, 1458679),
    @(13, -762, -121, -216, -1057, -1698),
    @(14, 128053, 211931, 194260, 766200, 1456981),
    @(15, "-", "-", 146398, 237894, 442130),
    @(16, "-", "-", -309320, -676270, -1467135),
    @(18, 121097, 112799, 2265, 0, 7700),
    @(19, -144547, -650422, -1187880, -1637924, -914267),
    @(20, 0, 0, 0, 0, 0),
    @(21, 0, 0, 0, 0, 0),
    @(22, -70, -850, -52984, -19610, -99),
    @(23, 0, 0, 0, 0, 0),
    @(24, 0, 0, 0, 0, 0),
    @(25, 0, 0, 0, 50000, 0),
    @(26, -2231, 0, -50000, 0, 0),
    @(27, 0, 0, 0, 0, 0),
    @(28, 0, 0, 0, 0, 0),
    @(29, 0, 0, 0, 0, 0),
    @(30, 0, 0, 0, 0, 0),
    @(31, 0, 0, 0, 0, 0),
    @(32, 0, 0, 0, 0, 0),
    @(33, 1742, 108, 1375, 6201, 1236),
    @(34, -24009, -538365, -1450146, -2039709, -1930435),
    @(35, 104044, -326434, -1255886, -1273509, -473454),
    @(37, 0, 0, 0, 147408, 0),
    @(38, "-", "-", 0, 0, 0),
    @(39, 0, 0, 0, 0, 0),
    @(40, 0, 0, 0, -7330, 0),
    @(41, 49000, 425000, 1318850, 1505182, 1258694),
    @(42, -59124, -2687, 0, 0, -538940),
    @(43, -9760, -61, 0, 0, -106291),
    @(44, 0, 0, 0, 0, 0),
    @(45, 0, 0, 0, 0, 0),
    @(46, 0, 0, 0, 0, 0),
    @(47, 0, 0, 0, 0, 0),
    @(48, 0, 0, 0, 0, 0),
    @(49, 0, 0, 0, 0, 0),
    @(50, 0, 0, 0, 0, 0),
    @(51, 0, 0, 0, 0, 0),
    @(52, -69630, -10384, -47234, -275753, -163969),
    @(53, -89514, 411868, 1271616, 1369507, 449494),
    @(54, 14530, 85434, 15730, 95998, -23960),
    @(55, 1192, 15722, 101177, 117028, 213034),
    @(56, 0, 21, 121, 8, 168),
    @(57, 15722, 101177, 117028, 213034, 189242),
    @(58, 0, 0, 0, 312592, 7047)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
}
